# Loan RBI, Variable Instalments
# The "Repayment schedule" sheet gets a new (blank) column inserted right
# before the "Late" column (column N, the 14th column) to make room for an
# additional variable-instalment figure. Everything from that column
# onward (Late / heading / Outstanding) shifts one column to the right.
# The sheet also becomes the active sheet/selection of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet (tabSelected moves here from
# "NewLoanInput", and the workbook's activeTab points at it too).
$ws.Activate()

# Insert a new blank column at position N (14), pushing the existing
# "Late" / heading / "Outstanding" columns one to the right.
$ws.Columns.Item(14).EntireColumn.Insert()

# The inherited width from the insert isn't right; match it to the
# neighbouring "In Advance" column (M) width, as in the authored file.
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Update the saved selection to the new last data column (now column R).
$ws.Range("R8").Select()
